$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 0")

$ws.Range("A308").Value = "a gremlin"
$ws.Range("B308").Value = "0x1B9C6"
$ws.Range("C308").Value = 2941
$ws.Range("D308").Value = 334

$ws.Range("A300").Value = "a goblin"
$ws.Range("B300").Value = "0x1B9B4"
$ws.Range("C300").Value = 0
$ws.Range("D300").Value = 723

$ws.Range("A15").Value = "an aegis mongbat"
$ws.Range("B15").Value = "0x1A11C"
$ws.Range("C15").Value = 2819
$ws.Range("D15").Value = 39

$ws.Range("A215").Value = "an effigy"
$ws.Range("B215").Value = "0x1A0E0"
$ws.Range("C215").Value = 2910
$ws.Range("D215").Value = 31

$ws.Range("A17").Value = "an aegis rat"
$ws.Range("B17").Value = "0x1A05C"
$ws.Range("C17").Value = 2944
$ws.Range("D17").Value = 215

$ws.Range("A8").Value = "an aegis cultist"
$ws.Range("B8").Value = "0x2432BE"
$ws.Range("C8").Value = 2928
$ws.Range("D8").Value = 689

$ws.Range("A227").Value = "an entombed"
$ws.Range("B227").Value = "0x2432A3"
$ws.Range("C227").Value = 2797
$ws.Range("D227").Value = 154

$ws.Range("A19").Value = "an aegis slime"
$ws.Range("B19").Value = "0x1B9CC"
$ws.Range("C19").Value = 2944
$ws.Range("D19").Value = 51

$ws.Range("A729").Value = "a wolfhound"
$ws.Range("B729").Value = "0x1A0E6"
$ws.Range("C729").Value = 2796
$ws.Range("D729").Value = 739

$ws.Range("A354").Value = "a kobold"
$ws.Range("B354").Value = "0x1A052"
$ws.Range("C354").Value = 0
$ws.Range("D354").Value = 724

$ws.Range("A123").Value = "a chaos footman"
$ws.Range("B123").Value = "0x1B9C4"
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 23

$ws.Range("A10").Value = "an aegis imp"
$ws.Range("B10").Value = "0x1B9D3"
$ws.Range("C10").Value = 2853
$ws.Range("D10").Value = 74

$ws.Range("A66").Value = "a blood ape"
$ws.Range("B66").Value = "0x1A0F6"
$ws.Range("C66").Value = 2466
$ws.Range("D66").Value = 29

$ws.Range("A14").Value = "an aegis minion"
$ws.Range("B14").Value = "0x1A079"
$ws.Range("C14").Value = 2928
$ws.Range("D14").Value = 776

$ws.Range("A82").Value = "a blood hellion"
$ws.Range("B82").Value = "0x1B9FD"
$ws.Range("C82").Value = 2740
$ws.Range("D82").Value = 241

$ws.Range("A86").Value = "a blood orc"
$ws.Range("B86").Value = "0x1B9EC"
$ws.Range("C86").Value = 2740
$ws.Range("D86").Value = 17

$ws.Range("A20").Value = "an aegis whelp"
$ws.Range("B20").Value = "0x18645"
$ws.Range("C20").Value = 2853
$ws.Range("D20").Value = 716

$ws.Range("A85").Value = "a blood ogre mage"
$ws.Range("B85").Value = "0x12AB9"
$ws.Range("C85").Value = 2846
$ws.Range("D85").Value = 1

$ws.Range("A91").Value = "a blood troll"
$ws.Range("B91").Value = "0x12AAD"
$ws.Range("C91").Value = 2846
$ws.Range("D91").Value = 53

$ws.Range("A130").Value = "a clay man"
$ws.Range("B130").Value = "0x12AAF"
$ws.Range("C130").Value = 2798
$ws.Range("D130").Value = 779

$ws.Range("A81").Value = "a blood harpy"
$ws.Range("B81").Value = "0xC360"
$ws.Range("C81").Value = 2740
$ws.Range("D81").Value = 30

$ws.Range("A77").Value = "a blood drake"
$ws.Range("B77").Value = "0x10303"
$ws.Range("C77").Value = 2740
$ws.Range("D77").Value = 61

$ws.Range("A16").Value = "an aegis noble"
$ws.Range("B16").Value = "0x10305"
$ws.Range("C16").Value = 1048
$ws.Range("D16").Value = 400

$ws.Range("A7").Value = "an aegis asp"
$ws.Range("B7").Value = "0x10356"
$ws.Range("C7").Value = 2853
$ws.Range("D7").Value = 727

$ws.Range("A18").Value = "an aegis scorpion"
$ws.Range("B18").Value = "0xE58B"
$ws.Range("C18").Value = 2944
$ws.Range("D18").Value = 717

$ws.Range("A13").Value = "an aegis lich"
$ws.Range("B13").Value = "0xE58C"
$ws.Range("C13").Value = 2944
$ws.Range("D13").Value = 24

$ws.Range("A92").Value = "a bloodrat"
$ws.Range("B92").Value = "0x18647"
$ws.Range("C92").Value = 2466
$ws.Range("D92").Value = 42

$ws.Range("A131").Value = "a coagulator"
$ws.Range("B131").Value = "0x12ABB"
$ws.Range("C131").Value = 2740
$ws.Range("D131").Value = 775

$ws.Range("A12").Value = "an aegis leech"
$ws.Range("B12").Value = "0x1111A"
$ws.Range("C12").Value = 2944
$ws.Range("D12").Value = 732

$ws.Range("A682").Value = "a twilight guardian"
$ws.Range("B682").Value = "0x1112B"
$ws.Range("C682").Value = 2274
$ws.Range("D682").Value = 199

$ws.Range("A601").Value = "a skulker"
$ws.Range("B601").Value = "0x12ABC"
$ws.Range("C601").Value = 2911
$ws.Range("D601").Value = 302

$ws.Range("A93").Value = "a bloodwolf"
$ws.Range("B93").Value = "0x107F9"
$ws.Range("C93").Value = 2740
$ws.Range("D93").Value = 1069

$ws.Range("A125").Value = "a chaos warrior"
$ws.Range("B125").Value = "0x255984"
$ws.Range("C125").Value = 0
$ws.Range("D125").Value = 25

$ws.Range("A74").Value = "a blood cyclops"
$ws.Range("B74").Value = "0x25596A"
$ws.Range("C74").Value = 2740
$ws.Range("D74").Value = 75

$ws.Range("A11").Value = "an aegis knight"
$ws.Range("B11").Value = "0x25583F"
$ws.Range("C11").Value = 1055
$ws.Range("D11").Value = 400

$ws.Range("A385").Value = "a malform"
$ws.Range("B385").Value = "0x1BA19"
$ws.Range("C385").Value = 2846
$ws.Range("D385").Value = 259

$ws.Range("A88").Value = "a blood scorpion"
$ws.Range("B88").Value = "0x1BA18"
$ws.Range("C88").Value = 2466
$ws.Range("D88").Value = 48

$ws.Range("A124").Value = "a chaos knight"
$ws.Range("B124").Value = "0x1BA24"
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 27

$ws.Range("A89").Value = "a blood serpent"
$ws.Range("B89").Value = "0x16B22"
$ws.Range("C89").Value = 2846
$ws.Range("D89").Value = 21

$ws.Range("A185").Value = "a doppelganger"
$ws.Range("B185").Value = "0x163E0"
$ws.Range("C185").Value = 2479
$ws.Range("D185").Value = 777

$ws.Range("A87").Value = "a blood ravager"
$ws.Range("B87").Value = "0x18649"
$ws.Range("C87").Value = 2740
$ws.Range("D87").Value = 314

$ws.Range("A78").Value = "blood elemental"
$ws.Range("B78").Value = "0xF78C"
$ws.Range("C78").Value = 2846
$ws.Range("D78").Value = 16

$ws.Range("A80").Value = "a blood fiend"
$ws.Range("B80").Value = "0x103BF"
$ws.Range("C80").Value = 2740
$ws.Range("D80").Value = 306

$ws.Range("A359").Value = "a lesser blood daemon"
$ws.Range("B359").Value = "0x103F7"
$ws.Range("C359").Value = 2740
$ws.Range("D359").Value = 40

$ws.Range("A94").Value = "a bloodworm"
$ws.Range("B94").Value = "0x103FD"
$ws.Range("C94").Value = 2817
$ws.Range("D94").Value = 287

$ws.Range("A76").Value = "a blood dragon"
$ws.Range("C76").Value = 2740
$ws.Range("D76").Value = 59

$ws.Range("A95").Value = "a blood dragon"
$ws.Range("B95").Value = "0x17EFF"
$ws.Range("C95").Value = 2740
$ws.Range("D95").Value = 59

$ws.Range("A79").Value = "a blood feaster"
$ws.Range("B79").Value = "0x101A5"
$ws.Range("C79").Value = 2740
$ws.Range("D79").Value = 303

$ws.Range("A83").Value = "a blood hunter"
$ws.Range("B83").Value = "0x17EF9"
$ws.Range("C83").Value = 2817
$ws.Range("D83").Value = 730

$ws.Range("A90").Value = "a blood sorcerer"
$ws.Range("B90").Value = "0x17F09"
$ws.Range("C90").Value = 2740
$ws.Range("D90").Value = 721

$ws.Range("A75").Value = "a blood daemon"
$ws.Range("B75").Value = "0x1BDD6"
$ws.Range("C75").Value = 2740
$ws.Range("D75").Value = 9

$ws.Range("A566").Value = "Sanguineous"
$ws.Range("B566").Value = "0x100B0"
$ws.Range("C566").Value = 2740
$ws.Range("D566").Value = 741

$ws.Range("B76").Value = "0x1B8A9"

$ws.Activate()
$ws.Range("D76").Select()
